$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'54.542.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.12%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.290.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.39%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.55%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'496.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.62%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'127.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.31%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.54%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.99%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.287.13"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.29%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.0951"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.29%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("E11").Value = "'  +2.36%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  +3.31%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  -2.48%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'2.666.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.21%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'21.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.00%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'54.375.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.92%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'  +0.91%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.280.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.57%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'10.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.83%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "'  +3.93%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'6.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +6.24%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'302.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.55%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.77%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  -1.72%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'62.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.04%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.28%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  +2.41%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.151"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.31%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "'2.384.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.29%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'7.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.73%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'169.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.48%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = "'  +0.45%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.0₃0689"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.14%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'5.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.96%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.26%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.09%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.12%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'17.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.26%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'1.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.25%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.868"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.39%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  +4.10%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'35.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.86%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  +3.08%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "'  +2.62%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  +2.08%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'128.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +5.04%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'4.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.68%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0892"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +1.80%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.544"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.29%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'240.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.09%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = "'  +3.04%  "
$ws.Range("E51").Style = "Normal"

